$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells stay text-typed (they were stored as text originally)
$ws.Range("C2:E3").NumberFormat = "@"

# Swap the runs/balls/fours values between row 2 and row 3
$ws.Range("C2").Value = "1"
$ws.Range("D2").Value = "2"
$ws.Range("E2").Value = "0"

$ws.Range("C3").Value = "12"
$ws.Range("D3").Value = "19"
$ws.Range("E3").Value = "1"
